# Generate Report for Handoff
# Updates the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps for the
# 07334b19-74e8-48ea-82ce-3e80974e1c39.md file row (row 5) across the three sheets
# to reflect a freshly generated handoff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-08-27 04:40:31"

# --- zh-cn sheet: column H = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-08-27 04:40:27"

# --- de-de sheet: column H = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-08-27 04:40:31"
